$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data contents of row 4 and row 6 (the "Valideringsstatus"
# column C value "Ovaliderad" is identical in both rows, so it is left as-is).

# --- Save row 4's original values ---
$r4_A  = $ws.Range("A4").Value()
$r4_B  = $ws.Range("B4").Value()
$r4_D  = $ws.Range("D4").Value()
$r4_E  = $ws.Range("E4").Value()
$r4_F  = $ws.Range("F4").Value()
$r4_G  = $ws.Range("G4").Value()
$r4_H  = $ws.Range("H4").Value()
$r4_S  = $ws.Range("S4").Value()
$r4_AC = $ws.Range("AC4").Value()
$r4_AI = $ws.Range("AI4").Value()

# --- Save row 6's original values ---
$r6_A  = $ws.Range("A6").Value()
$r6_B  = $ws.Range("B6").Value()
$r6_D  = $ws.Range("D6").Value()
$r6_E  = $ws.Range("E6").Value()
$r6_F  = $ws.Range("F6").Value()
$r6_G  = $ws.Range("G6").Value()
$r6_H  = $ws.Range("H6").Value()
$r6_S  = $ws.Range("S6").Value()

# --- Write row 6's original values into row 4 ---
$ws.Range("A4").Value = $r6_A
$ws.Range("B4").Value = $r6_B
$ws.Range("D4").Value = $r6_D
$ws.Range("E4").Value = $r6_E
$ws.Range("F4").Value = $r6_F
$ws.Range("G4").Value = $r6_G
$ws.Range("H4").Value = $r6_H
$ws.Range("S4").Value = $r6_S
$ws.Range("AC4").ClearContents()
$ws.Range("AI4").ClearContents()

# --- Write row 4's original values into row 6 ---
$ws.Range("A6").Value = $r4_A
$ws.Range("B6").Value = $r4_B
$ws.Range("D6").Value = $r4_D
$ws.Range("E6").Value = $r4_E
$ws.Range("F6").Value = $r4_F
$ws.Range("G6").Value = $r4_G
$ws.Range("H6").Value = $r4_H
$ws.Range("S6").Value = $r4_S
$ws.Range("AC6").Value = $r4_AC
$ws.Range("AI6").Value = $r4_AI
